# Append run: 2025-09-04 06:25 JST
# Replaces the scraped job-listing rows on the "ランサーズ" sheet (rows 2-9)
# with a fresh batch of results and drops the old rows 10-28.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Remove ALL existing hyperlinks first (collection-level delete is the
#    only reliable way to drop them in this engine) - we will re-create
#    the ones we still need (F2:F9) further down with the new URLs.
# ---------------------------------------------------------------------
$ws.Hyperlinks.Delete()

# ---------------------------------------------------------------------
# 2) Drop the old rows 10-28 entirely (sheet shrinks to A1:H9).
# ---------------------------------------------------------------------
$ws.Range("A10:H28").EntireRow.Delete()

# ---------------------------------------------------------------------
# 3) Write the new data for rows 2-9.
# ---------------------------------------------------------------------
$newRows = @(
    @{ Row = 2;  A = "2025-09-04 06:25:22"; B = "【急募】LINEで買取査定のAIシステム構築をお手伝いください!"; C = "システム開発"; D = "20,000 円 ~ 50,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5386178"; G = 313; H = "🔥AI,Ai" },
    @{ Row = 3;  A = "2025-09-04 06:25:22"; B = "日本株・米国株ランキングメール自動配信システムの作成依頼。Pythonなど。"; C = "システム開発"; D = "5,000 円 ~ 10,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5386223"; G = 205; H = "🔥Python" },
    @{ Row = 4;  A = "2025-09-04 06:25:22"; B = "オンラインスロットのスクレイピングソフトの制作"; C = "システム開発"; D = "200,000 円 ~ 300,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5386440"; G = 48;  H = "◆スクレイピング" },
    @{ Row = 5;  A = "2025-09-04 06:25:22"; B = "【急募】RUBYからPHPへのリプレース仕様書作成依頼"; C = "システム開発"; D = "1,000,000 円 ~ 3,000,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5386592"; G = 40;  H = "○PHP" },
    @{ Row = 6;  A = "2025-09-04 06:25:22"; B = "【緊急】運営しているサイトに表示される詐欺広告の削除方法を教えてください"; C = "システム開発"; D = "5,000 円 ~ 10,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5386516"; G = 30;  H = "◇サイト" },
    @{ Row = 7;  A = "2025-09-04 06:25:22"; B = "限定公開 PR 限定公開の仕事"; C = "システム開発"; D = "500,000 円 ~ 1,000,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5385681"; G = 25;  H = $null },
    @{ Row = 8;  A = "2025-09-04 06:25:22"; B = "注目 PR 超初級・SE育成の技術研修 サブ講師"; C = "システム開発"; D = "500,000 円 ~ 1,000,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5385021"; G = 25;  H = $null },
    @{ Row = 9;  A = "2025-09-04 06:25:22"; B = "限定公開 限定公開の仕事"; C = "システム開発"; D = "50,000 円 ~ 100,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5386235"; G = 18;  H = $null }
)

foreach ($item in $newRows) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
    $ws.Cells.Item($r, 8).Value = $item.H
}

# ---------------------------------------------------------------------
# 4) Re-create the hyperlinks for F2:F9 (in row order, so relationship
#    ids come out rId1..rId8) and restore the "Hyperlink" cell style
#    that Hyperlinks.Add slightly perturbs.
# ---------------------------------------------------------------------
foreach ($item in $newRows) {
    $r = $item.Row
    $cell = $ws.Cells.Item($r, 6)
    $ws.Hyperlinks.Add($cell, $item.F)
    $cell.Style = "Hyperlink"
}

# ---------------------------------------------------------------------
# 5) Column width tweaks: B 52 -> 40, H 23 -> 12 (raw OOXML character
#    units). Excel's ColumnWidth property is offset from the raw stored
#    width by 5/6 (~0.8333) of a character, so compensate for that.
# ---------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 40 - 5/6
$ws.Columns.Item(8).ColumnWidth = 12 - 5/6

Write-Output "edit complete"
